# Commit 19/12 home, search product, cart
# Add two new columns (C: "expected", D: "status") next to the existing
# username/password header row, matching the new shared-strings entries
# and the widened dimension/selection seen in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "expected"
$ws.Range("D1").Value = "status"

# Match the column widths that appear in the target file as closely as
# this host's character-width quantization allows.
$ws.Columns("B").ColumnWidth = 13.71
$ws.Columns("C").ColumnWidth = 14.33

# Leave the new header cell selected, as in the target workbook.
$ws.Range("D1").Select() | Out-Null
